# Add a new "hideInContents" column (S) to the survey, section1 and
# section2 sheets, and mark every "note" type row as hidden-in-contents
# (TRUE), to support hiding notes on the contents screen.

$wb = $excel.ActiveWorkbook

# --- survey sheet ---
$ws = $wb.Worksheets.Item("survey")
$ws.Cells.Item(1, 19).Value = "hideInContents"
$ws.Cells.Item(6, 19).Value = $true
$ws.Cells.Item(8, 19).Value = $true
$ws.Cells.Item(11, 19).Value = $true
$ws.Cells.Item(13, 19).Value = $true
$ws.Cells.Item(16, 19).Value = $true

# --- section1 sheet ---
$ws = $wb.Worksheets.Item("section1")
$ws.Cells.Item(1, 19).Value = "hideInContents"
$ws.Cells.Item(2, 19).Value = $true

# --- section2 sheet ---
$ws = $wb.Worksheets.Item("section2")
$ws.Cells.Item(1, 19).Value = "hideInContents"
$ws.Cells.Item(2, 19).Value = $true

# Leave the workbook's view state (active sheet/selection) the way it
# would be after a user reviews each edited sheet, ending back on survey.
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Activate() | Out-Null
$wsSettings.Range("A4:XFD4").Select() | Out-Null

$wsSection1 = $wb.Worksheets.Item("section1")
$wsSection1.Activate() | Out-Null
$wsSection1.Range("S3").Select() | Out-Null

$wsSection2 = $wb.Worksheets.Item("section2")
$wsSection2.Activate() | Out-Null
$wsSection2.Range("T2").Select() | Out-Null

$wsSurvey = $wb.Worksheets.Item("survey")
$wsSurvey.Activate() | Out-Null
$wsSurvey.Range("S16").Select() | Out-Null
